# netCrypto.xlsx update — relabel the "Roobic" payment type as "Wiretransfer"
# and bump the two USD Amount entries from hundreds to hundred-thousands.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared string "Roobic" becomes "Wiretransfer". Both cells that used
# that string (Payment Type N3 and InternalComment P3) need the new text.
$ws.Range("N3").Value = "Wiretransfer"
$ws.Range("P3").Value = "Wiretransfer"

# USD Amount column: 100 -> 100000, 50 -> 50000
$ws.Range("T2").Value = 100000
$ws.Range("T3").Value = 50000

# Move the active selection to T4 (matches the saved view state in the
# workbook after the edit).
$ws.Range("T4").Select()
